$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force column D to be treated as text so numeric-looking values
# (e.g. "0.999", "408.12") are stored as strings, matching the source data.
$dRange = $ws.Range('D2:D51')
$dRange.NumberFormat = '@'

# Row 2
$ws.Range('D2').Value = '61.764.79'
$ws.Range('E2').Value = '  +0.58%  '

# Row 3
$ws.Range('D3').Value = '3.419.08'
$ws.Range('E3').Value = '  +0.87%  '

# Row 4
$ws.Range('E4').Value = '  +0.13%  '

# Row 5
$ws.Range('D5').Value = '408.12'
$ws.Range('E5').Value = '  +0.97%  '

# Row 6
$ws.Range('D6').Value = '127.97'
$ws.Range('E6').Value = '  -1.71%  '

# Row 7
$ws.Range('D7').Value = '0.630'
$ws.Range('E7').Value = '  +7.12%  '

# Row 8
$ws.Range('E8').Value = '  -0.07%  '

# Row 9
$ws.Range('D9').Value = '0.732'
$ws.Range('E9').Value = '  +7.65%  '

# Row 10
$ws.Range('D10').Value = '0.138'
$ws.Range('E10').Value = '  +8.81%  '

# Row 11
$ws.Range('D11').Value = '42.52'
$ws.Range('E11').Value = '  +2.59%  '

# Row 12
$ws.Range('D12').Value = '9.11'
$ws.Range('E12').Value = '  +9.57%  '

# Row 13
$ws.Range('E13').Value = '  +0.19%  '

# Row 14
$ws.Range('B14').Value = 'Chainlink'
$ws.Range('C14').Value = 'https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link'
$ws.Range('D14').Value = '21.32'
$ws.Range('E14').Value = '  +8.38%  '

# Row 15
$ws.Range('B15').Value = 'WrappedliquidstakedEther2.0'
$ws.Range('C15').Value = 'https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth'
$ws.Range('D15').Value = '3.958.41'
$ws.Range('E15').Value = '  +0.61%  '

# Row 16
$ws.Range('D16').Value = '0.0000200'
$ws.Range('E16').Value = '  +41.80%  '

# Row 17
$ws.Range('D17').Value = '3.395.91'
$ws.Range('E17').Value = '  -0.48%  '

# Row 18
$ws.Range('D18').Value = '12.17'
$ws.Range('E18').Value = '  +4.62%  '

# Row 19
$ws.Range('E19').Value = '  +6.60%  '

# Row 20
$ws.Range('D20').Value = '61.736.75'
$ws.Range('E20').Value = '  +0.47%  '

# Row 21
$ws.Range('D21').Value = '440.26'
$ws.Range('E21').Value = '  +42.06%  '

# Row 22
$ws.Range('D22').Value = '91.32'
$ws.Range('E22').Value = '  +10.16%  '

# Row 23
$ws.Range('D23').Value = '3.18'
$ws.Range('E23').Value = '  +1.05%  '

# Row 24
$ws.Range('E24').Value = '  +2.28%  '

# Row 25
$ws.Range('D25').Value = '3.23'
$ws.Range('E25').Value = '  +2.96%  '

# Row 26
$ws.Range('D26').Value = '32.82'
$ws.Range('E26').Value = '  +11.89%  '

# Row 27
$ws.Range('D27').Value = '8.64'
$ws.Range('E27').Value = '  +7.86%  '

# Row 28
$ws.Range('E28').Value = '  +0.17%  '

# Row 29
$ws.Range('D29').Value = '7.63'
$ws.Range('E29').Value = '  -6.66%  '

# Row 30
$ws.Range('D30').Value = '2.71'
$ws.Range('E30').Value = '  +0.41%  '

# Row 31
$ws.Range('D31').Value = '11.93'
$ws.Range('E31').Value = '  +6.08%  '

# Row 32
$ws.Range('E32').Value = '  -0.28%  '

# Row 33
$ws.Range('D33').Value = '0.114'
$ws.Range('E33').Value = '  +0.26%  '

# Row 34
$ws.Range('D34').Value = '42.56'
$ws.Range('E34').Value = '  -2.71%  '

# Row 35
$ws.Range('D35').Value = '0.999'
$ws.Range('E35').Value = '  -0.05%  '

# Row 36
$ws.Range('D36').Value = '0.0497'
$ws.Range('E36').Value = '  +3.75%  '

# Row 37
$ws.Range('D37').Value = '53.23'
$ws.Range('E37').Value = '  +3.84%  '

# Row 38
$ws.Range('E38').Value = '  -0.10%  '

# Row 39
$ws.Range('D39').Value = '3.36'
$ws.Range('E39').Value = '  +0.71%  '

# Row 40
$ws.Range('E40').Value = '  +8.02%  '

# Row 41
$ws.Range('D41').Value = '2.93'
$ws.Range('E41').Value = '  -0.26%  '

# Row 42
$ws.Range('D42').Value = '142.50'
$ws.Range('E42').Value = '  +2.07%  '

# Row 43
$ws.Range('D43').Value = '0.311'
$ws.Range('E43').Value = '  -2.07%  '

# Row 44
$ws.Range('D44').Value = '4.20'
$ws.Range('E44').Value = '  +7.60%  '

# Row 45
$ws.Range('D45').Value = '1.98'
$ws.Range('E45').Value = '  +1.74%  '

# Row 46
$ws.Range('E46').Value = '  +13.55%  '

# Row 47
$ws.Range('D47').Value = '16.54'
$ws.Range('E47').Value = '  -0.36%  '

# Row 48
$ws.Range('D48').Value = '22.34'
$ws.Range('E48').Value = '  +6.68%  '

# Row 49
$ws.Range('E49').Value = '  +13.55%  '

# Row 50
$ws.Range('D50').Value = '3.769.45'
$ws.Range('E50').Value = '  +0.89%  '

# Row 51
$ws.Range('D51').Value = '2.118.45'
$ws.Range('E51').Value = '  +1.47%  '

# Restore default number format/style on column D so no stray styling remains
$dRange.NumberFormat = 'General'
$dRange.Style = 'Normal'
